# AoC "2023" sheet update — day 16 (row 17) stats corrected, day 17 (row 18) filled in.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2023")
$ws.Activate()

# --- Column B (column 2) and C (column 3) updates for rows 2-17 -----------
$updates = @{
    2  = @{ B = 214467; C = 68264 }
    3  = @{ B = 181154; C = 8390 }
    4  = @{ B = 119803; C = 17983 }
    5  = @{ B = 119182; C = 16002 }
    6  = @{ B = 72911;  C = 28205 }
    7  = @{ B = 93439;  C = 1550 }
    8  = @{ B = 73517;  C = 6477 }
    9  = @{ B = 66006;  C = 12991 }
    10 = @{ B = 67273;  C = 989 }
    11 = @{ B = 41784;  C = 15136 }
    12 = @{ B = 49389;  C = 2050 }
    13 = @{ B = 25252;  C = 13087 }
    14 = @{ B = 31306;  C = 4423 }
    15 = @{ B = 29157;  C = 6197 }
    16 = @{ B = 32013;  C = 3471 }
    17 = @{ B = 24889;  C = 915 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals.B
    $ws.Cells.Item($row, 3).Value = $vals.C
}

# --- Row 18 (AoC day 17) — previously blank, now filled in ----------------
$ws.Cells.Item(18, 2).Value = 8085   # B18
$ws.Cells.Item(18, 3).Value = 1093   # C18
$ws.Cells.Item(18, 5).Value = 7106   # E18
$ws.Cells.Item(18, 6).Value = 7793   # F18

# Row 18's formula cells were still carrying the cached "blank row" result
# (empty-string) from before B18/C18/E18/F18 had values, so re-assert the
# formulas in place to force a fresh evaluation against the new inputs.
$ws.Cells.Item(18, 4).Formula  = '=IF(ISBLANK(B18),"",B18+C18)'        # D18
$ws.Cells.Item(18, 7).Formula  = '=IF(D18="","",E18/D18)'              # G18
$ws.Cells.Item(18, 8).Formula  = '=IF(ISBLANK(C18),"",F18/B18)'        # H18
$ws.Cells.Item(18, 9).Formula  = '=IF(ISBLANK(E18),"",E18/$D$2)'       # I18
$ws.Cells.Item(18, 10).Formula = '=IF(ISBLANK(F18),"",F18/$B$2)'       # J18

# --- Selection moves from F18 to G18 ---------------------------------------
$ws.Range("G18").Select()
